$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2026-02-23 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-24 Tuesday", 2) | Out-Null

# Update the division problems in the table, addressed by (row, column)
# to disambiguate the two cells that originally shared the text "67÷3="
$t = $d.Tables(1)
$t.Cell(1, 1).Range.Text = "76÷4="
$t.Cell(1, 2).Range.Text = "41÷4="
$t.Cell(1, 3).Range.Text = "57÷3="
$t.Cell(1, 4).Range.Text = "43÷2="
$t.Cell(1, 5).Range.Text = "35÷4="
$t.Cell(5, 1).Range.Text = "81÷3="
$t.Cell(5, 2).Range.Text = "88÷8="
$t.Cell(5, 3).Range.Text = "99÷6="
$t.Cell(5, 4).Range.Text = "13÷9="
$t.Cell(5, 5).Range.Text = "61÷6="
$t.Cell(9, 1).Range.Text = "74÷7="
$t.Cell(9, 2).Range.Text = "88÷4="
$t.Cell(9, 3).Range.Text = "19÷3="
$t.Cell(9, 4).Range.Text = "24÷5="
$t.Cell(9, 5).Range.Text = "15÷3="
$t.Cell(13, 1).Range.Text = "83÷7="
$t.Cell(13, 2).Range.Text = "20÷4="
$t.Cell(13, 3).Range.Text = "30÷7="
$t.Cell(13, 4).Range.Text = "79÷8="
$t.Cell(13, 5).Range.Text = "16÷6="
$t.Cell(17, 1).Range.Text = "20÷2="
$t.Cell(17, 2).Range.Text = "35÷8="
$t.Cell(17, 3).Range.Text = "22÷5="
$t.Cell(17, 4).Range.Text = "59÷5="
$t.Cell(17, 5).Range.Text = "61÷4="
